# Scheduled-runner style refresh of the per-sheet price/profit columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) with freshly pulled
# market-board data. Only numeric value cells change; labels/formatting
# are untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H49").Value = 1084
$ws.Range("I49").Value = 228.4
$ws.Range("J49").Value = 3223
$ws.Range("K49").Value = 685.2
$ws.Range("L49").Value = 9669
$ws.Range("M49").Value = -549.2
$ws.Range("N49").Value = -9941

$ws.Range("H62").Value = 30854.854
$ws.Range("I62").Value = 38993.1
$ws.Range("K62").Value = 38993.1
$ws.Range("M62").Value = -38369.1

$ws.Range("H64").Value = 13105.125
$ws.Range("I64").Value = 10132.333
$ws.Range("K64").Value = 10132.333
$ws.Range("M64").Value = -9884.333000000001

$ws.Range("H65").Value = 30854.854
$ws.Range("I65").Value = 38993.1
$ws.Range("K65").Value = 194965.5
$ws.Range("M65").Value = -191845.5

$ws.Range("H67").Value = 13105.125
$ws.Range("I67").Value = 10132.333
$ws.Range("K67").Value = 10132.333
$ws.Range("M67").Value = -9274.333000000001

$ws.Range("H94").Value = 29604.533
$ws.Range("I94").Value = 31654.5
$ws.Range("J94").Value = 905
$ws.Range("K94").Value = 31654.5
$ws.Range("L94").Value = 905
$ws.Range("M94").Value = -31203.5
$ws.Range("N94").Value = -1807

$ws.Range("H112").Value = 87381.086
$ws.Range("J112").Value = 95290.73
$ws.Range("L112").Value = 285872.19
$ws.Range("N112").Value = -288088.19

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5536.121
$ws.Range("I45").Value = 8813.200000000001
$ws.Range("J45").Value = 2805.2222
$ws.Range("K45").Value = 8813.200000000001
$ws.Range("L45").Value = 2805.2222
$ws.Range("M45").Value = -8436.200000000001
$ws.Range("N45").Value = -3559.2222

$ws.Range("H51").Value = 9642
$ws.Range("I51").Value = 9642
$ws.Range("K51").Value = 9642
$ws.Range("M51").Value = -8886

$ws.Range("H56").Value = 17500
$ws.Range("J56").Value = 25000
$ws.Range("L56").Value = 25000
$ws.Range("N56").Value = -26484

$ws.Range("H61").Value = 3451.0417
$ws.Range("I61").Value = 3240.9473
$ws.Range("K61").Value = 3240.9473
$ws.Range("M61").Value = -3028.9473

$ws.Range("H74").Value = 7526.9287
$ws.Range("I74").Value = 879.55
$ws.Range("K74").Value = 879.55
$ws.Range("M74").Value = -5.549999999999955

$ws.Range("H77").Value = 7526.9287
$ws.Range("I77").Value = 879.55
$ws.Range("K77").Value = 4397.75
$ws.Range("M77").Value = -29.75

$ws.Range("H122").Value = 1562.9412
$ws.Range("I122").Value = 1076.5714
$ws.Range("K122").Value = 3229.7142
$ws.Range("M122").Value = -779.7142000000003

$ws.Range("H132").Value = 2460.158
$ws.Range("I132").Value = 2061.9666
$ws.Range("K132").Value = 6185.899800000001
$ws.Range("M132").Value = -3655.899800000001

$ws.Range("H136").Value = 3451.0417
$ws.Range("I136").Value = 3240.9473
$ws.Range("K136").Value = 9722.841899999999
$ws.Range("M136").Value = -7172.841899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5036.25
$ws.Range("I20").Value = 4947.8335
$ws.Range("K20").Value = 4947.8335
$ws.Range("M20").Value = -4700.8335

$ws.Range("H22").Value = 346931.12
$ws.Range("J22").Value = 530757
$ws.Range("L22").Value = 530757
$ws.Range("N22").Value = -531103

$ws.Range("H134").Value = 2496.9783
$ws.Range("I134").Value = 2216.9487
$ws.Range("J134").Value = 4057.1428
$ws.Range("K134").Value = 6650.8461
$ws.Range("L134").Value = 12171.4284
$ws.Range("M134").Value = -4115.8461
$ws.Range("N134").Value = -17241.4284

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 33616
$ws.Range("J41").Value = 49999
$ws.Range("L41").Value = 49999
$ws.Range("N41").Value = -50855

$ws.Range("H50").Value = 10500
$ws.Range("I50").Value = 6500
$ws.Range("K50").Value = 6500
$ws.Range("M50").Value = -5875

$ws.Range("H51").Value = 90
$ws.Range("I51").Value = 90
$ws.Range("K51").Value = 90
$ws.Range("M51").Value = 646

$ws.Range("H59").Value = 25649.95
$ws.Range("J59").Value = 25684.21
$ws.Range("L59").Value = 25684.21
$ws.Range("N59").Value = -27974.21

$ws.Range("H60").Value = 7490.091
$ws.Range("I60").Value = 6239.2
$ws.Range("K60").Value = 6239.2
$ws.Range("M60").Value = -5728.2

$ws.Range("H61").Value = 90
$ws.Range("I61").Value = 90
$ws.Range("K61").Value = 90
$ws.Range("M61").Value = 258

$ws.Range("H74").Value = 40000
$ws.Range("J74").Value = 40000
$ws.Range("L74").Value = 40000
$ws.Range("N74").Value = -41748

$ws.Range("H77").Value = 40000
$ws.Range("J77").Value = 40000
$ws.Range("L77").Value = 120000
$ws.Range("N77").Value = -128736

$ws.Range("H134").Value = 8916.25
$ws.Range("I134").Value = 7394.1714
$ws.Range("J134").Value = 13014.154
$ws.Range("K134").Value = 22182.5142
$ws.Range("L134").Value = 39042.462
$ws.Range("M134").Value = -19647.5142
$ws.Range("N134").Value = -44112.462

$ws.Range("H135").Value = 74933.336
$ws.Range("J135").Value = 74933.336
$ws.Range("L135").Value = 74933.336
$ws.Range("N135").Value = -85073.336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 139.7
$ws.Range("I11").Value = 69.40000000000001
$ws.Range("J11").Value = 210
$ws.Range("K11").Value = 208.2
$ws.Range("L11").Value = 630
$ws.Range("M11").Value = -68.20000000000002
$ws.Range("N11").Value = -910

$ws.Range("H58").Value = 650
$ws.Range("I58").Value = 650
$ws.Range("K58").Value = 1950
$ws.Range("M58").Value = -1822

$ws.Range("H87").Value = 26355.4
$ws.Range("I87").Value = 20888.5
$ws.Range("J87").Value = 30000
$ws.Range("K87").Value = 62665.5
$ws.Range("L87").Value = 90000
$ws.Range("M87").Value = -61417.5
$ws.Range("N87").Value = -92496

$ws.Range("H90").Value = 26355.4
$ws.Range("I90").Value = 20888.5
$ws.Range("J90").Value = 30000
$ws.Range("K90").Value = 187996.5
$ws.Range("L90").Value = 270000
$ws.Range("M90").Value = -181756.5
$ws.Range("N90").Value = -282480

$ws.Range("H103").Value = 5658.3335
$ws.Range("I103").Value = 4950
$ws.Range("J103").Value = 5800
$ws.Range("K103").Value = 14850
$ws.Range("L103").Value = 17400
$ws.Range("M103").Value = -13971
$ws.Range("N103").Value = -19158

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12386.077
$ws.Range("I70").Value = 10117.8
$ws.Range("K70").Value = 10117.8
$ws.Range("M70").Value = -9847.799999999999

$ws.Range("H73").Value = 12386.077
$ws.Range("I73").Value = 10117.8
$ws.Range("K73").Value = 10117.8
$ws.Range("M73").Value = -9181.799999999999

$ws.Range("H80").Value = 6233.1816
$ws.Range("I80").Value = 3110.8333
$ws.Range("J80").Value = 9980
$ws.Range("K80").Value = 3110.8333
$ws.Range("L80").Value = 9980
$ws.Range("M80").Value = -2112.8333
$ws.Range("N80").Value = -11976

$ws.Range("H83").Value = 6233.1816
$ws.Range("I83").Value = 3110.8333
$ws.Range("J83").Value = 9980
$ws.Range("K83").Value = 15554.1665
$ws.Range("L83").Value = 49900
$ws.Range("M83").Value = -10562.1665
$ws.Range("N83").Value = -59884

$ws.Range("H97").Value = 1482
$ws.Range("I97").Value = 1552.5
$ws.Range("K97").Value = 1552.5
$ws.Range("M97").Value = -1056.5

$ws.Range("H122").Value = 1171.0834
$ws.Range("I122").Value = 1116.1111
$ws.Range("J122").Value = 1336
$ws.Range("K122").Value = 3348.3333
$ws.Range("L122").Value = 4008
$ws.Range("M122").Value = -898.3333000000002
$ws.Range("N122").Value = -8908

$ws.Range("H126").Value = 10490.902
$ws.Range("I126").Value = 14084.667
$ws.Range("K126").Value = 42254.001
$ws.Range("M126").Value = -39784.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 915.5
$ws.Range("I100").Value = 915.5
$ws.Range("K100").Value = 915.5
$ws.Range("M100").Value = -374.5

$ws.Range("H136").Value = 6083.5264
$ws.Range("I136").Value = 5276.154
$ws.Range("K136").Value = 15828.462
$ws.Range("M136").Value = -13278.462

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 10000
$ws.Range("J37").Value = 10000
$ws.Range("L37").Value = 10000
$ws.Range("N37").Value = -10406

$ws.Range("H81").Value = 8216.727999999999
$ws.Range("I81").Value = 16274.286
$ws.Range("J81").Value = 4456.533
$ws.Range("K81").Value = 32548.572
$ws.Range("L81").Value = 8913.066000000001
$ws.Range("M81").Value = -31487.572
$ws.Range("N81").Value = -11035.066

$ws.Range("H84").Value = 8216.727999999999
$ws.Range("I84").Value = 16274.286
$ws.Range("J84").Value = 4456.533
$ws.Range("K84").Value = 162742.86
$ws.Range("L84").Value = 44565.33
$ws.Range("M84").Value = -157438.86
$ws.Range("N84").Value = -55173.33
